# Auto-generated edit script: updates the cryptocurrency price/volume table
# on the active worksheet to match the 'Updated cryptos list' GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. '53.50', '1.010',
# '0.00000000358') need an explicit Text number format, otherwise Excel
# auto-converts the assigned string into a numeric value and mangles the
# formatting (trailing zeros, scientific notation, etc). After the value is
# written we restore the cell style to 'Normal' so no extra formatting is
# left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '29.678.36'
$ws.Range("E2").Value = '  -1.38%  '
# Row 3
$ws.Range("D3").Value = '2.097.00'
$ws.Range("E3").Value = '  -0.40%  '
# Row 4
$ws.Range("E4").Value = '  +0.47%  '
# Row 5
Set-TextValue $ws.Range("D5") '343.27'
$ws.Range("E5").Value = '  -1.99%  '
# Row 6
$ws.Range("E6").Value = '  +0.49%  '
# Row 7
$ws.Range("E7").Value = '  +0.34%  '
# Row 8
Set-TextValue $ws.Range("D8") '0.4378'
$ws.Range("E8").Value = '  -1.90%  '
# Row 9
Set-TextValue $ws.Range("D9") '53.50'
$ws.Range("E9").Value = '  +1.85%  '
# Row 10
Set-TextValue $ws.Range("D10") '0.09198'
# Row 11
$ws.Range("E11").Value = '  -1.08%  '
# Row 12
Set-TextValue $ws.Range("D12") '24.62'
$ws.Range("E12").Value = '  -4.91%  '
# Row 13
Set-TextValue $ws.Range("D13") '6.767'
$ws.Range("E13").Value = '  +0.25%  '
# Row 14
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D14") '8.140'
$ws.Range("E14").Value = '  -1.32%  '
# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.035.01'
$ws.Range("E15").Value = '  -3.41%  '
# Row 16
Set-TextValue $ws.Range("D16") '102.57'
$ws.Range("E16").Value = '  +3.42%  '
# Row 17
$ws.Range("E17").Value = '  -0.02%  '
# Row 18
Set-TextValue $ws.Range("D18") '1.010'
$ws.Range("E18").Value = '  +0.50%  '
# Row 19
$ws.Range("E19").Value = '  +1.35%  '
# Row 20
Set-TextValue $ws.Range("D20") '0.06664'
$ws.Range("E20").Value = '  -0.15%  '
# Row 21
Set-TextValue $ws.Range("D21") '1.007'
$ws.Range("E21").Value = '  +0.46%  '
# Row 22
Set-TextValue $ws.Range("D22") '6.207'
$ws.Range("E22").Value = '  -0.83%  '
# Row 23
$ws.Range("D23").Value = '29.740.55'
$ws.Range("E23").Value = '  -1.46%  '
# Row 24
$ws.Range("E24").Value = '  -2.39%  '
# Row 25
Set-TextValue $ws.Range("D25") '2.304'
$ws.Range("E25").Value = '  -1.75%  '
# Row 26
$ws.Range("D26").Value = '2.328.61'
$ws.Range("E26").Value = '  -1.21%  '
# Row 27
$ws.Range("E27").Value = '  -0.80%  '
# Row 28
Set-TextValue $ws.Range("D28") '161.84'
$ws.Range("E28").Value = '  -0.46%  '
# Row 29
Set-TextValue $ws.Range("D29") '2.493'
$ws.Range("E29").Value = '  -2.69%  '
# Row 30
Set-TextValue $ws.Range("D30") '133.51'
$ws.Range("E30").Value = '  -0.16%  '
# Row 31
$ws.Range("E31").Value = '  -4.67%  '
# Row 32
Set-TextValue $ws.Range("D32") '1.685'
$ws.Range("E32").Value = '  +2.77%  '
# Row 33
$ws.Range("E33").Value = '  -1.80%  '
# Row 34
$ws.Range("E34").Value = '  -1.35%  '
# Row 35
Set-TextValue $ws.Range("D35") '3.959'
$ws.Range("E35").Value = '  -0.33%  '
# Row 36
Set-TextValue $ws.Range("D36") '6.396'
$ws.Range("E36").Value = '  +7.78%  '
# Row 37
Set-TextValue $ws.Range("D37") '10.41'
$ws.Range("E37").Value = '  +0.70%  '
# Row 38
Set-TextValue $ws.Range("D38") '0.02577'
$ws.Range("E38").Value = '  -0.31%  '
# Row 39
Set-TextValue $ws.Range("D39") '0.06710'
$ws.Range("E39").Value = '  -2.14%  '
# Row 40
Set-TextValue $ws.Range("D40") '0.6974'
$ws.Range("E40").Value = '  +2.01%  '
# Row 41
Set-TextValue $ws.Range("D41") '1.337'
$ws.Range("E41").Value = '  +6.32%  '
# Row 42
Set-TextValue $ws.Range("D42") '12.44'
$ws.Range("E42").Value = '  -2.19%  '
# Row 43
Set-TextValue $ws.Range("D43") '0.2214'
$ws.Range("E43").Value = '  -4.32%  '
# Row 44
Set-TextValue $ws.Range("D44") '0.6767'
$ws.Range("E44").Value = '  +5.40%  '
# Row 45
Set-TextValue $ws.Range("D45") '14.32'
$ws.Range("E45").Value = '  -0.50%  '
# Row 46
$ws.Range("E46").Value = '  -0.59%  '
# Row 47
Set-TextValue $ws.Range("D47") '0.00000000358'
$ws.Range("E47").Value = '  -1.23%  '
# Row 48
Set-TextValue $ws.Range("D48") '3.622'
$ws.Range("E48").Value = '  -1.07%  '
# Row 49
Set-TextValue $ws.Range("D49") '1.202'
$ws.Range("E49").Value = '  +3.27%  '
# Row 50
$ws.Range("E50").Value = '  -0.69%  '
# Row 51
Set-TextValue $ws.Range("D51") '81.05'
$ws.Range("E51").Value = '  -2.80%  '
